$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.171978831291199
$ws.Range("B1").Value = 2.436966419219971
$ws.Range("D1").Value = 2.365468263626099
$ws.Range("E1").Value = 1.237959742546082
